$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.719.22"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.599.08"
$ws.Range("E3").Value = "  +0.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue "D5" "211.57"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.10%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.24%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +0.76%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.91%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.823.67"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.604.36"
$ws.Range("E13").Value = "  -0.28%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.89%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.48%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.04"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0₃0740"
$ws.Range("E17").Value = "  -3.77%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "208.73"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.14"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.44%  "

# Row 22 - Toncoin
$ws.Range("E22").Value = "  -3.85%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.02"
$ws.Range("E23").Value = "  +1.00%  "

# Row 24 - Monero
Set-TextValue "D24" "143.56"
$ws.Range("E24").Value = "  +0.46%  "

# Row 25 - BinanceUSD
$ws.Range("E25").Value = "  +0.08%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +0.23%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.47%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.33"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29 - Hedera
$ws.Range("E29").Value = "  -2.18%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.05%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +0.66%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.40%  "

# Row 33/34 - Maker and WEMIXToken swap places
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D33" "1.25"
$ws.Range("E33").Value = "  +17.53%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D34" "1.275.98"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.56%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -4.27%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -1.15%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +0.11%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  +0.54%  "

# Row 41 - MXToken
$ws.Range("E41").Value = "  +0.27%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  -0.59%  "

# Row 43 - Aave
Set-TextValue "D43" "62.53"
$ws.Range("E43").Value = "  -1.01%  "

# Row 44 - RocketPoolETH
Set-TextValue "D44" "1.735.15"
$ws.Range("E44").Value = "  +0.14%  "

# Row 45 - Quant
Set-TextValue "D45" "90.42"
$ws.Range("E45").Value = "  -0.68%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +0.75%  "

# Row 47 - Algorand
$ws.Range("E47").Value = "  +1.39%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +0.80%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "7.57"
$ws.Range("E49").Value = "  +3.81%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  +1.57%  "
